$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fix row 7, column A style (currently default, should match A2:A6 bordered style)
$ws.Range("A6").Copy()
$ws.Range("A7:A11").PasteSpecial(-4122)  # xlPasteFormats

# Copy formatting for new rows 8-11, columns B and D, from row 7 (already correctly styled)
$ws.Range("B7").Copy()
$ws.Range("B8:B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D7").Copy()
$ws.Range("D8:D11").PasteSpecial(-4122)  # xlPasteFormats

# Copy formatting for new rows 8-11, column C, from C7 (style s=7)
$ws.Range("C7").Copy()
$ws.Range("C8:C11").PasteSpecial(-4122)  # xlPasteFormats

# Update Row 7 Results: PASS -> SKIP
$ws.Range("D7").Value = "SKIP"

# Fill in new rows in the order the author appears to have entered them
$ws.Range("B8").Value = "To verify that document count gets decreased in the watchlist page when a document is deleted from watchlist"
$ws.Range("B9").Value = "To verify that MORE button doesn't get displayed if number of documents in watchlist page is less than or equal to 10"

$ws.Range("A8").Value = "WatchlistDeleteArticleTest"
$ws.Range("A9").Value = "WatchlistMoreButtonBelowTenArticlesTest"
$ws.Range("A10").Value = "WatchlistMoreButtonAboveTenArticlesTest"
$ws.Range("A11").Value = "NavigateToWatchlistFromRVTest"

$ws.Range("B11").Value = "To verify that app navigates to correct page when user navigates back from document page"
$ws.Range("D11").Value = "FAIL"

$ws.Range("B10").Value = "To verify that MORE button is present in watchlist page if total search results is more than 10`nTo verify that MORE button is working correctly in watchlist page`n"

$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("C10").Value = "Y"
$ws.Range("C11").Value = "Y"

$ws.Range("D8").Value = "SKIP"
$ws.Range("D9").Value = "SKIP"
$ws.Range("D10").Value = "SKIP"

# Row heights
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 61.5

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B12").Select()
